$wb = $excel.ActiveWorkbook

# Update the "№ студенческого билета, зачетки" header to "N_ZACHET" on every sheet (C3)
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("C3").Value = "N_ZACHET"
}

# Sheet "134" (4th sheet): remove the blank row 14 that separated the
# regular roster rows from the two manually-appended rows, shifting the
# "Алоян Артём" / "Сахно Антон" rows up from 15/16 to 14/15.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(14).Delete()

# Update the selection on each non-active sheet first …
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws2.Range("C3").Select()
$ws3.Range("C3").Select()
$ws4.Range("B28").Select()

# … then make sheet "131" the active tab with C3 selected, matching the
# saved view state of the edited workbook.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("C3").Select()
